# Revise build script optimization
# Adds newly captured sensor readings to the FE-LIFTER and MID-LIFTER sheets.

$G_CONST = [double]"5.686312626471138e+23"

function Add-Row {
    param($ws, $rowNum, $a, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Range("A$rowNum").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A$rowNum").Value2 = $a
    $ws.Range("B$rowNum").Value2 = $b
    $ws.Range("C$rowNum").Value2 = $c
    $ws.Range("D$rowNum").Value2 = $d
    $ws.Range("E$rowNum").Value2 = $e
    $ws.Range("F$rowNum").Value2 = $f
    $ws.Range("G$rowNum").Value2 = $g
    $ws.Range("H$rowNum").Value2 = $h
    $ws.Range("I$rowNum").Value2 = $i
}

$wb = $excel.ActiveWorkbook

# ROW50-FE-LIFTER (sheet 1): add rows 26-28
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-Row $ws1 26 45729.08020857639 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20
Add-Row $ws1 27 45729.08023003472 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20
Add-Row $ws1 28 45729.08025329861 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20

# ROW50-MID-LIFTER (sheet 2): add rows 59-67
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-Row $ws2 59 45729.06475774306 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 60 45729.06477990741 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 61 45729.06480305555 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 62 45729.14823328704 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 63 45729.14825528935 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 64 45729.14827864584 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 65 45729.23170871528 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 66 45729.23173081018 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws2 67 45729.23175395833 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25

# ROW11-FE-LIFTER (sheet 3): add rows 26-28
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-Row $ws3 26 45729.08020857639 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20
Add-Row $ws3 27 45729.08023003472 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20
Add-Row $ws3 28 45729.08025329861 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 $G_CONST 400 20

# ROW11-MID-LIFTER (sheet 4): add rows 59-67
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-Row $ws4 59 45729.06475774306 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 60 45729.06477990741 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 61 45729.06480305555 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 62 45729.14823328704 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 63 45729.14825528935 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 64 45729.14827864584 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 65 45729.23170871528 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 66 45729.23173081018 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25
Add-Row $ws4 67 45729.23175395833 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 $G_CONST 400 25

Write-Host "Edit complete"
